$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 179.1580256666667
$ws.Range("H2").Value = 537.4740770000001
$ws.Range("I2").Value = 0.3468013736386751
$ws.Range("J2").Value = 0.3468013736386751
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.237840333333334
$ws.Range("N2").Value = 12.713521
$ws.Range("O2").Value = 0.4475167411500002
$ws.Range("P2").Value = 0.4475167411500002
$ws.Range("Q2").Value = 759.2431072105687
$ws.Range("R2").Value = 6833.187964895118
$ws.Range("S2").Value = 0.1551994205571235
$ws.Range("T2").Value = 0.1551994205571235

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 179.1580256666667
$ws.Range("H3").Value = 537.4740770000001
$ws.Range("I3").Value = 0.3468013736386751
$ws.Range("J3").Value = 0.3468013736386751
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.231839666666667
$ws.Range("N3").Value = 15.695519
$ws.Range("O3").Value = 0.5524832588499998
$ws.Range("P3").Value = 0.5524832588499998
$ws.Range("Q3").Value = 937.3260652845516
$ws.Range("R3").Value = 8435.934587560963
$ws.Range("S3").Value = 0.1916019530815517
$ws.Range("T3").Value = 0.1916019530815516

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 239.807332
$ws.Range("H4").Value = 719.421996
$ws.Range("I4").Value = 0.4642019905988459
$ws.Range("J4").Value = 0.4642019905988459
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.237840333333334
$ws.Range("N4").Value = 12.713521
$ws.Range("O4").Value = 0.4475167411500002
$ws.Range("P4").Value = 0.4475167411500002
$ws.Range("Q4").Value = 1016.265183778657
$ws.Range("R4").Value = 9146.386654007916
$ws.Range("S4").Value = 0.2077381620681386
$ws.Range("T4").Value = 0.2077381620681386

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 239.807332
$ws.Range("H5").Value = 719.421996
$ws.Range("I5").Value = 0.4642019905988459
$ws.Range("J5").Value = 0.4642019905988459
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.231839666666667
$ws.Range("N5").Value = 15.695519
$ws.Range("O5").Value = 0.5524832588499998
$ws.Range("P5").Value = 0.5524832588499998
$ws.Range("Q5").Value = 1254.633511915103
$ws.Range("R5").Value = 11291.70160723592
$ws.Range("S5").Value = 0.2564638285307074
$ws.Range("T5").Value = 0.2564638285307074

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 97.63589966666666
$ws.Range("H6").Value = 292.907699
$ws.Range("I6").Value = 0.1889966357624789
$ws.Range("J6").Value = 0.1889966357624789
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 4.237840333333334
$ws.Range("N6").Value = 12.713521
$ws.Range("O6").Value = 0.4475167411500002
$ws.Range("P6").Value = 0.4475167411500002
$ws.Range("Q6").Value = 413.7653535886865
$ws.Range("R6").Value = 3723.888182298179
$ws.Range("S6").Value = 0.08457915852473816
$ws.Range("T6").Value = 0.08457915852473816

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 97.63589966666666
$ws.Range("H7").Value = 292.907699
$ws.Range("I7").Value = 0.1889966357624789
$ws.Range("J7").Value = 0.1889966357624789
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.231839666666667
$ws.Range("N7").Value = 15.695519
$ws.Range("O7").Value = 0.5524832588499998
$ws.Range("P7").Value = 0.5524832588499998
$ws.Range("Q7").Value = 510.8153727667534
$ws.Range("R7").Value = 4597.33835490078
$ws.Range("S7").Value = 0.1044174772377408
$ws.Range("T7").Value = 0.1044174772377408

